# The deck's slide master theme ("Integral") and the notes-master theme
# ("Office Theme") had their colour schemes swapped: the slide master now
# carries the stock "Office" palette. Reproduce that by rewriting each of
# the 12 theme colour slots on the (single) Design's slide-master theme
# through ThemeColorScheme, the supported colour-scheme edit surface.

function ToRGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p  = $ppt.ActivePresentation
$d  = $p.Designs.Item(1)
$sm = $d.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

# Target palette = the stock Office theme colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink), matching ThemeColorScheme.Colors(1..12).
$tcs.Colors(1).RGB  = ToRGB 0x00 0x00 0x00   # dk1      000000
$tcs.Colors(2).RGB  = ToRGB 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Colors(3).RGB  = ToRGB 0x44 0x54 0x6A   # dk2      44546A
$tcs.Colors(4).RGB  = ToRGB 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Colors(5).RGB  = ToRGB 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Colors(6).RGB  = ToRGB 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Colors(7).RGB  = ToRGB 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Colors(8).RGB  = ToRGB 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Colors(9).RGB  = ToRGB 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Colors(10).RGB = ToRGB 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Colors(11).RGB = ToRGB 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Colors(12).RGB = ToRGB 0x95 0x4F 0x72   # folHlink 954F72
